$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = -8.429
$ws.Range("D10").Value = -8.231999999999999
$ws.Range("D12").Value = -7.290000000000001
$ws.Range("D18").Value = -8.318999999999999
$ws.Range("D25").Value = -8.228999999999999
$ws.Range("D37").Value = -8.175000000000001
$ws.Range("D55").Value = -8.318999999999999
$ws.Range("D68").Value = -7.111
$ws.Range("D77").Value = -8.198
$ws.Range("D78").Value = -8.406000000000001
$ws.Range("D79").Value = -7.85
$ws.Range("D80").Value = -7.923
$ws.Range("D81").Value = -7.274000000000001
$ws.Range("D82").Value = -8.399000000000001
$ws.Range("D84").Value = -8.315000000000001
$ws.Range("D101").Value = -7.767
$ws.Range("D102").Value = -8.081999999999999
